$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-05-27 Tuesday" "2025-05-28 Wednesday"

Replace-Text "155×7=" "306×2="
Replace-Text "760×5=" "285×2="
Replace-Text "563×6=" "808×7="
Replace-Text "315×9=" "581×8="
Replace-Text "437×2=" "711×7="
Replace-Text "855×4=" "728×5="
Replace-Text "261×7=" "225×9="
Replace-Text "127×7=" "743×4="
Replace-Text "426×6=" "952×5="
Replace-Text "395×9=" "196×2="
Replace-Text "954×3=" "115×5="
Replace-Text "823×4=" "564×7="
Replace-Text "212×7=" "290×6="
Replace-Text "167×7=" "573×8="
Replace-Text "921×3=" "688×3="
Replace-Text "422×6=" "261×9="
Replace-Text "302×7=" "213×7="
Replace-Text "454×9=" "171×7="
Replace-Text "124×9=" "578×8="
Replace-Text "957×9=" "147×7="
Replace-Text "349×7=" "847×8="
Replace-Text "523×6=" "706×4="
Replace-Text "653×2=" "268×2="
Replace-Text "972×3=" "521×6="
Replace-Text "923×2=" "394×4="
